# Adds an extra blank paragraph immediately before each of the four
# "incident.*" narrative paragraphs in the domestic Statement of Facts
# template (emotional effect / physical injury / witnesses / police
# actions), giving the incident section additional spacing.
#
# For each target paragraph we locate it with Find, then call
# InsertParagraphAfter() on its *previous* paragraph's range -- this
# creates the new empty paragraph using the previous paragraph's
# paragraph-mark formatting (ind w:left="465" w:right="577", no
# firstLine) while leaving the target paragraph (and its firstLine
# indent) completely untouched, exactly as in the source edit.

$d = $word.ActiveDocument

$anchors = @(
    "Because of this incident I felt",
    "As a result, I was physically injured",
    "There were witnesses to the incident",
    "The police were called and did the following"
)

foreach ($anchor in $anchors) {
    $rng = $d.Content
    $rng.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $targetPara = $rng.Paragraphs(1)
    $prevPara = $targetPara.Previous()
    $prevPara.Range.InsertParagraphAfter()
}
